$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing value for week 50 (row 51)
$ws.Range("B51").Value = 564

# Add new data rows for weeks 51 and 52
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 406

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 36
